$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) column stores plain text (e.g. "18.50", "214.47") that would
# otherwise be auto-coerced into a Double by Excel (losing trailing zeros /
# introducing floating-point noise). Force those cells to Text format first.
$priceCells = @("D2", "D3", "D5", "D6", "D10", "D12", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D30", "D31", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D44", "D47", "D48", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.091.47'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.640.50'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.63%  '
$ws.Range("D5").Value = '214.47'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = '0.505'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("E8").Value = '  -2.51%  '
$ws.Range("E9").Value = '  -2.32%  '
$ws.Range("D10").Value = '18.50'
$ws.Range("E10").Value = '  -5.78%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '1.647.71'
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("E14").Value = '  -2.82%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '62.30'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0₃0749'
$ws.Range("E16").Value = '  -1.97%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.089.29'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").Value = '190.14'
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("D20").Value = '4.27'
$ws.Range("E20").Value = '  -2.37%  '
$ws.Range("D21").Value = '9.53'
$ws.Range("E21").Value = '  -3.98%  '
$ws.Range("D22").Value = '6.12'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").Value = '144.19'
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = '0.0485'
$ws.Range("E30").Value = '  -3.64%  '
$ws.Range("D31").Value = '3.18'
$ws.Range("E31").Value = '  -2.03%  '
$ws.Range("E32").Value = '  -3.89%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").Value = '1.51'
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("D36").Value = '1.124.14'
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").Value = '2.46'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '0.523'
$ws.Range("E38").Value = '  -4.17%  '
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").Value = '98.84'
$ws.Range("E40").Value = '  -0.51%  '
$ws.Range("D41").Value = '0.785'
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("D42").Value = '5.28'
$ws.Range("E42").Value = '  -3.68%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '55.18'
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").Value = '0.415'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").Value = '7.61'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").Value = '0.0931'
$ws.Range("E50").Value = '  -3.42%  '
$ws.Range("E51").Value = '  -1.19%  '

# Restore the default (Normal) cell style so no stray formatting beyond the
# content itself is introduced.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
